$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Multiple Employers")

# Row 5 (D5:G5): round the daily rate to the nearest integer (assuming 100% sickness)
$fmt = '_(* #,##0_);_(* \(#,##0\);_(* "-"??_);_(@_)'
$ws.Range("D5:G5").NumberFormat = $fmt

$ws.Range("D5").Formula = "=ROUND(D3*12/260,0)"
$ws.Range("E5:G5").Formula = "=ROUND(E3*12/260,0)"

# H5 note changes from "Double" to "Integer" now that the value is rounded
$ws.Range("H5").Value = "Integer"

# New note in I5 explaining the change, highlighted in red
$ws.Range("I5").Value = "New: Rounding assuming 100% sickness"
$ws.Range("I5").Font.Color = 255

# Reflect the active selection used while making the edit
$ws.Range("I6").Select()
